$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header-independent data cells (rows 2-51) to match target state.
# Column A = rank/id (numeric), B = Nazwa, C = Poziom, D = Wojewodztwo, E = Skutecznosc[%]

# Row 2
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "ty"
$ws.Cells.Item(2, 3).Value = "Extreme"
$ws.Cells.Item(2, 4).Value = "Opolskie"
$ws.Cells.Item(2, 5).Value = 100

# Row 3
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "testowańsko"
$ws.Cells.Item(3, 3).Value = "Extreme"
$ws.Cells.Item(3, 4).Value = "Opolskie"
$ws.Cells.Item(3, 5).Value = 100

# Row 4
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "j"
$ws.Cells.Item(4, 3).Value = "Extreme"
$ws.Cells.Item(4, 4).Value = "Opolskie"
$ws.Cells.Item(4, 5).Value = 100

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "k"
$ws.Cells.Item(5, 3).Value = "Extreme"
$ws.Cells.Item(5, 4).Value = "Opolskie"
$ws.Cells.Item(5, 5).Value = 92

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Formula = "=""'"""
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(6, 2).PasteSpecial(-4163)
$ws.Cells.Item(6, 3).Value = "Extreme"
$ws.Cells.Item(6, 4).Value = "Opolskie"
$ws.Cells.Item(6, 5).Value = 92

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "testowańsko"
$ws.Cells.Item(7, 3).Value = "Extreme"
$ws.Cells.Item(7, 4).Value = "Opolskie"
$ws.Cells.Item(7, 5).Value = 92

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "11111"
$ws.Cells.Item(8, 3).Value = "Extreme"
$ws.Cells.Item(8, 4).Value = "Lubelskie"
$ws.Cells.Item(8, 5).Value = 88

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Maks"
$ws.Cells.Item(9, 3).Value = "Extreme"
$ws.Cells.Item(9, 4).Value = "Dolnośląskie"
$ws.Cells.Item(9, 5).Value = 87

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "maks"
$ws.Cells.Item(10, 3).Value = "Extreme"
$ws.Cells.Item(10, 4).Value = "Wszystkie"
$ws.Cells.Item(10, 5).Value = 85

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "k"
$ws.Cells.Item(11, 3).Value = "Extreme"
$ws.Cells.Item(11, 4).Value = "Opolskie"
$ws.Cells.Item(11, 5).Value = 83

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Maks"
$ws.Cells.Item(12, 3).Value = "Extreme"
$ws.Cells.Item(12, 4).Value = "Lubelskie"
$ws.Cells.Item(12, 5).Value = 72

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "l"
$ws.Cells.Item(13, 3).Value = "Extreme"
$ws.Cells.Item(13, 4).Value = "Opolskie"
$ws.Cells.Item(13, 5).Value = 67

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "mak"
$ws.Cells.Item(14, 3).Value = "Easy"
$ws.Cells.Item(14, 4).Value = "Opolskie"
$ws.Cells.Item(14, 5).Value = 67

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "pkstz"
$ws.Cells.Item(15, 3).Value = "Extreme"
$ws.Cells.Item(15, 4).Value = "Zachodniopomorskie"
$ws.Cells.Item(15, 5).Value = 66

# Row 16
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "wiki <3"
$ws.Cells.Item(16, 3).Value = "Extreme"
$ws.Cells.Item(16, 4).Value = "Lubelskie"
$ws.Cells.Item(16, 5).Value = 60

# Row 17
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "testyyyyyyy"
$ws.Cells.Item(17, 3).Value = "Extreme"
$ws.Cells.Item(17, 4).Value = "Kujawsko-Pomorskie"
$ws.Cells.Item(17, 5).Value = 57

# Row 18
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "k"
$ws.Cells.Item(18, 3).Value = "Hard"
$ws.Cells.Item(18, 4).Value = "Lubuskie"
$ws.Cells.Item(18, 5).Value = 50

# Row 19
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Maks"
$ws.Cells.Item(19, 3).Value = "Extreme"
$ws.Cells.Item(19, 4).Value = "Podlaskie"
$ws.Cells.Item(19, 5).Value = 48

# Row 20
$ws.Cells.Item(20, 1).Value = 20
$ws.Cells.Item(20, 2).Value = "a"
$ws.Cells.Item(20, 3).Value = "Medium"
$ws.Cells.Item(20, 4).Value = "Opolskie"
$ws.Cells.Item(20, 5).Value = 42

# Row 21
$ws.Cells.Item(21, 1).Value = 21
$ws.Cells.Item(21, 2).Value = "Karolcio"
$ws.Cells.Item(21, 3).Value = "Extreme"
$ws.Cells.Item(21, 4).Value = "Wielkopolskie"
$ws.Cells.Item(21, 5).Value = 42

# Row 22
$ws.Cells.Item(22, 1).Value = 19
$ws.Cells.Item(22, 2).Value = "a"
$ws.Cells.Item(22, 3).Value = "Extreme"
$ws.Cells.Item(22, 4).Value = "Lubelskie"
$ws.Cells.Item(22, 5).Value = 42

# Row 23
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "maks"
$ws.Cells.Item(23, 3).Value = "Extreme"
$ws.Cells.Item(23, 4).Value = "Wszystkie"
$ws.Cells.Item(23, 5).Value = 40

# Row 24
$ws.Cells.Item(24, 1).Value = 24
$ws.Cells.Item(24, 2).Value = "a"
$ws.Cells.Item(24, 3).Value = "Hard"
$ws.Cells.Item(24, 4).Value = "Lubuskie"
$ws.Cells.Item(24, 5).Value = 36

# Row 25
$ws.Cells.Item(25, 1).Value = 25
$ws.Cells.Item(25, 2).ClearContents()
$ws.Cells.Item(25, 3).Value = "Extreme"
$ws.Cells.Item(25, 4).Value = "Dolnośląskie"
$ws.Cells.Item(25, 5).Value = 36

# Row 26
$ws.Cells.Item(26, 1).Value = 23
$ws.Cells.Item(26, 2).Value = "Maks"
$ws.Cells.Item(26, 3).Value = "Extreme"
$ws.Cells.Item(26, 4).Value = "Opolskie"
$ws.Cells.Item(26, 5).Value = 36

# Row 27
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "kkk"
$ws.Cells.Item(27, 3).Value = "Easy"
$ws.Cells.Item(27, 4).Value = "Podlaskie"
$ws.Cells.Item(27, 5).Value = 35

# Row 28
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "Maks"
$ws.Cells.Item(28, 3).Value = "Medium"
$ws.Cells.Item(28, 4).Value = "Podlaskie"
$ws.Cells.Item(28, 5).Value = 35

# Row 29
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).ClearContents()
$ws.Cells.Item(29, 3).Value = "Extreme"
$ws.Cells.Item(29, 4).Value = "Podlaskie"
$ws.Cells.Item(29, 5).Value = 35

# Row 30
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "szymek"
$ws.Cells.Item(30, 3).Value = "Extreme"
$ws.Cells.Item(30, 4).Value = "Kujawsko-Pomorskie"
$ws.Cells.Item(30, 5).Value = 35

# Row 31
$ws.Cells.Item(31, 1).Value = 34
$ws.Cells.Item(31, 2).Value = "h"
$ws.Cells.Item(31, 3).Value = "Hard"
$ws.Cells.Item(31, 4).Value = "Małopolskie"
$ws.Cells.Item(31, 5).Value = 33

# Row 32
$ws.Cells.Item(32, 1).Value = 37
$ws.Cells.Item(32, 2).Value = "maks"
$ws.Cells.Item(32, 3).Value = "Extreme"
$ws.Cells.Item(32, 4).Value = "Śląskie"
$ws.Cells.Item(32, 5).Value = 33

# Row 33
$ws.Cells.Item(33, 1).Value = 36
$ws.Cells.Item(33, 2).Value = "hjk"
$ws.Cells.Item(33, 3).Value = "Extreme"
$ws.Cells.Item(33, 4).Value = "Lubelskie"
$ws.Cells.Item(33, 5).Value = 33

# Row 34
$ws.Cells.Item(34, 1).Value = 35
$ws.Cells.Item(34, 2).ClearContents()
$ws.Cells.Item(34, 3).Value = "Extreme"
$ws.Cells.Item(34, 4).Value = "Opolskie"
$ws.Cells.Item(34, 5).Value = 33

# Row 35
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = "h"
$ws.Cells.Item(35, 3).Value = "Extreme"
$ws.Cells.Item(35, 4).Value = "Opolskie"
$ws.Cells.Item(35, 5).Value = 33

# Row 36
$ws.Cells.Item(36, 1).Value = 32
$ws.Cells.Item(36, 2).Value = "pkstz"
$ws.Cells.Item(36, 3).Value = "Extreme"
$ws.Cells.Item(36, 4).Value = "Zachodniopomorskie"
$ws.Cells.Item(36, 5).Value = 33

# Row 37
$ws.Cells.Item(37, 1).Value = 31
$ws.Cells.Item(37, 2).Value = "uj"
$ws.Cells.Item(37, 3).Value = "Medium"
$ws.Cells.Item(37, 4).Value = "Opolskie"
$ws.Cells.Item(37, 5).Value = 33

# Row 38
$ws.Cells.Item(38, 1).Value = 30
$ws.Cells.Item(38, 2).Value = "maks"
$ws.Cells.Item(38, 3).Value = "Extreme"
$ws.Cells.Item(38, 4).Value = "Śląskie"
$ws.Cells.Item(38, 5).Value = 33

# Row 39
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Formula = "=""'"""
$ws.Cells.Item(39, 2).Copy()
$ws.Cells.Item(39, 2).PasteSpecial(-4163)
$ws.Cells.Item(39, 3).Value = "Extreme"
$ws.Cells.Item(39, 4).Value = "Łódzkie"
$ws.Cells.Item(39, 5).Value = 32

# Row 40
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "aa"
$ws.Cells.Item(40, 3).Value = "Extreme"
$ws.Cells.Item(40, 4).Value = "Dolnośląskie"
$ws.Cells.Item(40, 5).Value = 32

# Row 41
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).ClearContents()
$ws.Cells.Item(41, 3).Value = "Extreme"
$ws.Cells.Item(41, 4).Value = "Dolnośląskie"
$ws.Cells.Item(41, 5).Value = 30

# Row 42
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).ClearContents()
$ws.Cells.Item(42, 3).Value = "Medium"
$ws.Cells.Item(42, 4).Value = "Lubuskie"
$ws.Cells.Item(42, 5).Value = 29

# Row 43
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "hjk"
$ws.Cells.Item(43, 3).Value = "Extreme"
$ws.Cells.Item(43, 4).Value = "Lubelskie"
$ws.Cells.Item(43, 5).Value = 29

# Row 44
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "Maks"
$ws.Cells.Item(44, 3).Value = "Extreme"
$ws.Cells.Item(44, 4).Value = "Wszystkie"
$ws.Cells.Item(44, 5).Value = 28

# Row 45
$ws.Cells.Item(45, 1).Value = 0
$ws.Cells.Item(45, 2).Value = "iu"
$ws.Cells.Item(45, 3).Value = "Extreme"
$ws.Cells.Item(45, 4).Value = "Łódzkie"
$ws.Cells.Item(45, 5).Value = 28

# Row 46
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = "Maks"
$ws.Cells.Item(46, 3).Value = "Extreme"
$ws.Cells.Item(46, 4).Value = "Zachodniopomorskie"
$ws.Cells.Item(46, 5).Value = 27

# Row 47
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).ClearContents()
$ws.Cells.Item(47, 3).Value = "Extreme"
$ws.Cells.Item(47, 4).Value = "Łódzkie"
$ws.Cells.Item(47, 5).Value = 27

# Row 48
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = "Maks"
$ws.Cells.Item(48, 3).Value = "Extreme"
$ws.Cells.Item(48, 4).Value = "Dolnośląskie"
$ws.Cells.Item(48, 5).Value = 27

# Row 49
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = "a"
$ws.Cells.Item(49, 3).Value = "Extreme"
$ws.Cells.Item(49, 4).Value = "Dolnośląskie"
$ws.Cells.Item(49, 5).Value = 26

# Row 50
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = "mak"
$ws.Cells.Item(50, 3).Value = "Extreme"
$ws.Cells.Item(50, 4).Value = "Opolskie"
$ws.Cells.Item(50, 5).Value = 25

# Row 51
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = ";"
$ws.Cells.Item(51, 3).Value = "Extreme"
$ws.Cells.Item(51, 4).Value = "Opolskie"
$ws.Cells.Item(51, 5).Value = 25

$excel.CutCopyMode = $false
Write-Host "Applied V17.0 multiplier fixes"
